# feat: add 2022-Q1 data
#
# The workbook currently has two sheets:
#   1. "2021-Q4" - per-fund holding detail for the 2021-Q4 quarter
#   2. "总计"     - summary (date / count / value) with a single 2021-Q4 row
#
# This script:
#   1. Duplicates the "2021-Q4" sheet (preserving all formatting/styles) and
#      places the copy right after it, renaming it to "2022-Q1".
#   2. Overwrites the copied sheet's data cells with the 2022-Q1 fund detail
#      values (numeric-looking text values are entered with a leading
#      apostrophe so they stay text cells, matching the source data).
#   3. Updates the "总计" summary sheet: the existing 2021-Q4 summary row
#      moves down to row 3, and a new row 2 is added for 2022-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: duplicate "2021-Q4" right after itself, then rename the copy.
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $q4)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q1"

# ---------------------------------------------------------------------
# Step 2: overwrite the copied sheet's data with the 2022-Q1 values.
# Row 2: fund 004044 (rank 3), Row 3: fund 002303 (rank 2).
# ---------------------------------------------------------------------

# Row 2 - 004044 / 金鹰转型动力灵活配置混合
# (fund code and the numeric-looking figures are entered with a leading
# apostrophe so Excel keeps them as text, matching the source data -
# otherwise e.g. "004044" would lose its leading zeros as a number.)
$newSheet.Range("B2").Formula = "'004044"
$newSheet.Range("C2").Value = "金鹰转型动力灵活配置混合"
$newSheet.Range("D2").Formula = "'0.72"
$newSheet.Range("E2").Formula = "'93.34"
$newSheet.Range("F2").Formula = "'7.78"
$newSheet.Range("G2").Formula = "'0.0560"
$newSheet.Range("H2").Value = 3

# Row 3 - 002303 / 金鹰智慧生活灵活配置混合
$newSheet.Range("B3").Formula = "'002303"
$newSheet.Range("C3").Value = "金鹰智慧生活灵活配置混合"
$newSheet.Range("D3").Formula = "'0.11"
$newSheet.Range("E3").Formula = "'89.88"
$newSheet.Range("F3").Formula = "'8.09"
$newSheet.Range("G3").Formula = "'0.0089"
$newSheet.Range("H3").Value = 2

# ---------------------------------------------------------------------
# Step 3: update the "总计" summary sheet with the new quarter.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Push the existing 2021-Q4 summary row from row 2 down to row 3,
# copying A2 first so the row-index style (s="2") comes along, then
# fixing up the value.
$total.Range("A2").Copy($total.Range("A3"))
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.22

# New row 2 for 2022-Q1.
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.06
